$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all cell contents (keeps the header rows style) so the shared-string table
# can be rebuilt from scratch in the exact order the cells are (re)written below -
# this mirrors the append order produced by the upstream script regeneration.
$ws.Cells.ClearContents()

# --- Re-write header row (row 1) ---
$ws.Cells.Item(1, 1).Value = "Sending cluster"
$ws.Cells.Item(1, 2).Value = "Ligand symbol"
$ws.Cells.Item(1, 3).Value = "Receptor symbol"
$ws.Cells.Item(1, 4).Value = "Target cluster"
$ws.Cells.Item(1, 5).Value = "Ligand-expressing cells"
$ws.Cells.Item(1, 6).Value = "Ligand detection rate"
$ws.Cells.Item(1, 7).Value = "Ligand average expression value"
$ws.Cells.Item(1, 8).Value = "Ligand total expression value"
$ws.Cells.Item(1, 9).Value = "Ligand derived specificity of average expression value"
$ws.Cells.Item(1, 10).Value = "Ligand derived specificity of total expression value"
$ws.Cells.Item(1, 11).Value = "Receptor-expressing cells"
$ws.Cells.Item(1, 12).Value = "Receptor detection rate"
$ws.Cells.Item(1, 13).Value = "Receptor average expression value"
$ws.Cells.Item(1, 14).Value = "Receptor total expression value"
$ws.Cells.Item(1, 15).Value = "Receptor derived specificity of average expression value"
$ws.Cells.Item(1, 16).Value = "Receptor derived specificity of total expression value"
$ws.Cells.Item(1, 17).Value = "Edge average expression weight"
$ws.Cells.Item(1, 18).Value = "Edge total expression weight"
$ws.Cells.Item(1, 19).Value = "Edge average expression derived specificity"
$ws.Cells.Item(1, 20).Value = "Edge total expression derived specificity"

# --- Write string columns A-D column-major so the shared-string table ends up
# ordered ECs, FAPs, Fgl1, Lag3, MuSCs, Resolving-Mac (indices 20-25) ---
# Column A
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(9, 1).Value = "FAPs"

# Column B
$ws.Cells.Item(2, 2).Value = "Fgl1"
$ws.Cells.Item(3, 2).Value = "Fgl1"
$ws.Cells.Item(4, 2).Value = "Fgl1"
$ws.Cells.Item(5, 2).Value = "Fgl1"
$ws.Cells.Item(6, 2).Value = "Fgl1"
$ws.Cells.Item(7, 2).Value = "Fgl1"
$ws.Cells.Item(8, 2).Value = "Fgl1"
$ws.Cells.Item(9, 2).Value = "Fgl1"

# Column C
$ws.Cells.Item(2, 3).Value = "Lag3"
$ws.Cells.Item(3, 3).Value = "Lag3"
$ws.Cells.Item(4, 3).Value = "Lag3"
$ws.Cells.Item(5, 3).Value = "Lag3"
$ws.Cells.Item(6, 3).Value = "Lag3"
$ws.Cells.Item(7, 3).Value = "Lag3"
$ws.Cells.Item(8, 3).Value = "Lag3"
$ws.Cells.Item(9, 3).Value = "Lag3"

# Column D
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"

# --- Update numeric columns (E-T) ---
# Row 2
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.004391333333333333
$ws.Cells.Item(2, 8).Value = 0.013174
$ws.Cells.Item(2, 9).Value = 0.04203587120571539
$ws.Cells.Item(2, 10).Value = 0.04203587120571539
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 10.718847
$ws.Cells.Item(2, 14).Value = 32.156541
$ws.Cells.Item(2, 15).Value = 0.2473266771098565
$ws.Cells.Item(2, 16).Value = 0.2473266771098565
$ws.Cells.Item(2, 17).Value = 0.04707003012599999
$ws.Cells.Item(2, 18).Value = 0.423630271134
$ws.Cells.Item(2, 19).Value = 0.01039659234472749
$ws.Cells.Item(2, 20).Value = 0.01039659234472749

# Row 3
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.004391333333333333
$ws.Cells.Item(3, 8).Value = 0.013174
$ws.Cells.Item(3, 9).Value = 0.04203587120571539
$ws.Cells.Item(3, 10).Value = 0.04203587120571539
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 10.56216766666667
$ws.Cells.Item(3, 14).Value = 31.686503
$ws.Cells.Item(3, 15).Value = 0.2437114581515935
$ws.Cells.Item(3, 16).Value = 0.2437114581515935
$ws.Cells.Item(3, 17).Value = 0.04638199894688888
$ws.Cells.Item(3, 18).Value = 0.417437990522
$ws.Cells.Item(3, 19).Value = 0.01024462346621748
$ws.Cells.Item(3, 20).Value = 0.01024462346621748

# Row 4
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.004391333333333333
$ws.Cells.Item(4, 8).Value = 0.013174
$ws.Cells.Item(4, 9).Value = 0.04203587120571539
$ws.Cells.Item(4, 10).Value = 0.04203587120571539
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 7.214691666666667
$ws.Cells.Item(4, 14).Value = 21.644075
$ws.Cells.Item(4, 15).Value = 0.1664717964804274
$ws.Cells.Item(4, 16).Value = 0.1664717964804274
$ws.Cells.Item(4, 17).Value = 0.03168211600555555
$ws.Cells.Item(4, 18).Value = 0.28513904405
$ws.Cells.Item(4, 19).Value = 0.006997786996235311
$ws.Cells.Item(4, 20).Value = 0.006997786996235311

# Row 5
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.004391333333333333
$ws.Cells.Item(5, 8).Value = 0.013174
$ws.Cells.Item(5, 9).Value = 0.04203587120571539
$ws.Cells.Item(5, 10).Value = 0.04203587120571539
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 14.84311633333333
$ws.Cells.Item(5, 14).Value = 44.529349
$ws.Cells.Item(5, 15).Value = 0.3424900682581226
$ws.Cells.Item(5, 16).Value = 0.3424900682581225
$ws.Cells.Item(5, 17).Value = 0.06518107152511111
$ws.Cells.Item(5, 18).Value = 0.586629643726
$ws.Cells.Item(5, 19).Value = 0.01439686839853511
$ws.Cells.Item(5, 20).Value = 0.01439686839853511

# Row 6
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.100075
$ws.Cells.Item(6, 8).Value = 0.300225
$ws.Cells.Item(6, 9).Value = 0.9579641287942846
$ws.Cells.Item(6, 10).Value = 0.9579641287942845
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 10.718847
$ws.Cells.Item(6, 14).Value = 32.156541
$ws.Cells.Item(6, 15).Value = 0.2473266771098565
$ws.Cells.Item(6, 16).Value = 0.2473266771098565
$ws.Cells.Item(6, 17).Value = 1.072688613525
$ws.Cells.Item(6, 18).Value = 9.654197521724999
$ws.Cells.Item(6, 19).Value = 0.2369300847651291
$ws.Cells.Item(6, 20).Value = 0.236930084765129

# Row 7
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.100075
$ws.Cells.Item(7, 8).Value = 0.300225
$ws.Cells.Item(7, 9).Value = 0.9579641287942846
$ws.Cells.Item(7, 10).Value = 0.9579641287942845
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 10.56216766666667
$ws.Cells.Item(7, 14).Value = 31.686503
$ws.Cells.Item(7, 15).Value = 0.2437114581515935
$ws.Cells.Item(7, 16).Value = 0.2437114581515935
$ws.Cells.Item(7, 17).Value = 1.057008929241667
$ws.Cells.Item(7, 18).Value = 9.513080363175
$ws.Cells.Item(7, 19).Value = 0.233466834685376
$ws.Cells.Item(7, 20).Value = 0.233466834685376

# Row 8
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.100075
$ws.Cells.Item(8, 8).Value = 0.300225
$ws.Cells.Item(8, 9).Value = 0.9579641287942846
$ws.Cells.Item(8, 10).Value = 0.9579641287942845
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 7.214691666666667
$ws.Cells.Item(8, 14).Value = 21.644075
$ws.Cells.Item(8, 15).Value = 0.1664717964804274
$ws.Cells.Item(8, 16).Value = 0.1664717964804274
$ws.Cells.Item(8, 17).Value = 0.7220102685416667
$ws.Cells.Item(8, 18).Value = 6.498092416875001
$ws.Cells.Item(8, 19).Value = 0.1594740094841921
$ws.Cells.Item(8, 20).Value = 0.1594740094841921

# Row 9
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.100075
$ws.Cells.Item(9, 8).Value = 0.300225
$ws.Cells.Item(9, 9).Value = 0.9579641287942846
$ws.Cells.Item(9, 10).Value = 0.9579641287942845
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 14.84311633333333
$ws.Cells.Item(9, 14).Value = 44.529349
$ws.Cells.Item(9, 15).Value = 0.3424900682581226
$ws.Cells.Item(9, 16).Value = 0.3424900682581225
$ws.Cells.Item(9, 17).Value = 1.485424867058334
$ws.Cells.Item(9, 18).Value = 13.368823803525
$ws.Cells.Item(9, 19).Value = 0.3280931998595875
$ws.Cells.Item(9, 20).Value = 0.3280931998595873
